$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. OrangeHRM_LoginForm (sheet 4): trim stale trailing column from the used
#    range (column F had no data, just leftover dimension/spans metadata).
# ---------------------------------------------------------------------------
$wsOrange = $wb.Worksheets.Item(4)
$null = $wsOrange.UsedRange

# ---------------------------------------------------------------------------
# 2. ParaBank_RegistartionForm (sheet 5): fix the spelling of the sheet name
#    and refresh its sample data with the new Alpha/Beta test users.
# ---------------------------------------------------------------------------
$wsRegistration = $wb.Worksheets.Item(5)
$wsRegistration.Name = "ParaBank_RegistrationForm"

$wsRegistration.Range("A2").Value = "Alpha"
$wsRegistration.Range("B2").Value = "Beta"
$wsRegistration.Range("I2").Value = "Alpha"
$wsRegistration.Range("M2").Value = "Your account was created successfully. You are now logged in."
$wsRegistration.Range("N2").Value = "FAIL"
$wsRegistration.Range("N2").Select()

# ---------------------------------------------------------------------------
# 3. Add a brand-new ParaBank_LoginForm sheet after the registration sheet,
#    covering a successful login and a failed login with the new users.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsLogin = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsLogin.Name = "ParaBank_LoginForm"

$wsLogin.Range("A1").Value = "Username"
$wsLogin.Range("B1").Value = "Password"
$wsLogin.Range("C1").Value = "Expected Results"
$wsLogin.Range("D1").Value = "Actual Results"
$wsLogin.Range("E1").Value = "Status"

$wsLogin.Range("A2").Value = "Alpha"
$wsLogin.Range("B2").Value = "secret123"
$wsLogin.Range("E2").Value = "PASS"

$wsLogin.Range("A3").Value = "Beta"
$wsLogin.Range("B3").Value = "wrongpassword"
$wsLogin.Range("C3").Value = "The username and password could not be verified."
$wsLogin.Range("D3").Value = "The username and password could not be verified."
$wsLogin.Range("E3").Value = "PASS"

$wsLogin.Range("C2").Value = "Login Success"
$wsLogin.Range("D2").Value = "Login Success"

$wsLogin.Columns.Item(1).ColumnWidth = 9
$wsLogin.Columns.Item(2).ColumnWidth = 14.5
$wsLogin.Columns.Item(3).ColumnWidth = 46
$wsLogin.Columns.Item(4).ColumnWidth = 20.5
$wsLogin.Columns.Item(5).ColumnWidth = 5.5

$wsLogin.Range("E2").Select()
